$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (Pont en H / LB1930MC-AH) : quantity set to 0, row hidden ---
$ws.Range("C3").Value = 0
$ws.Rows.Item(3).Hidden = $true

# --- Row 5 (CAN I2C / MAX1237EUA+) : quantity set to 0, row hidden ---
$ws.Range("C5").Value = 0
$ws.Rows.Item(5).Hidden = $true

# --- Row 7 (trigger de schmitt / SN74LVC3G17DCTR) : quantity set to 0, row hidden + shorter row height ---
$ws.Range("C7").Value = 0
$ws.Rows.Item(7).Hidden = $true
$ws.Rows.Item(7).RowHeight = 13.5

# --- Row 8 : new component "pont en H (v2,0)" / LB1836M-TLM-E ---
$ws.Range("E8").Value = "http://fr.rs-online.com/web/p/drivers-de-moteur/7570499/?searchTerm=LB1836M-TLM-E&relevancy-data=636F3D3226696E3D4931384E4B6E6F776E41734D504E266C753D6672266D6D3D6D61746368616C6C7061727469616C26706D3D5E5B5C707B4C7D5C707B4E647D2D2C2F255C2E5D2B2426706F3D313326736E3D592673743D4D414E5F504152545F4E554D4245522677633D424F5448267573743D4C42313833364D2D544C4D2D4526&sra=p"
$ws.Range("B8").Value = "pont en H (v2,0)"
$ws.Range("D8").Value = "LB1836M-TLM-E"
$ws.Range("C8").Value = 2
$ws.Range("F8").Value = 2.65

# Style the new reference cell (D8) like a pasted RS-Online product name: bold, small Arial, dark grey
$f8 = $ws.Range("D8").Font
$f8.Name = "Arial"
$f8.Size = 9
$f8.Bold = $true
$f8.Color = 3355443

# Hyperlink marker left on E8 (no live target, as in the source workbook)
$ws.Hyperlinks.Add($ws.Range("E8"), "") | Out-Null

# --- Row 10 : new component (SN74AHC14N) ---
$ws.Range("D10").Value = " SN74AHC14N"
$ws.Range("E10").Value = "http://www.ti.com/lit/ds/symlink/sn74ahc14.pdf"
$ws.Range("B10").Value = "trigger de schmitt inverseur"
$ws.Range("C10").Value = 1

$f10 = $ws.Range("D10").Font
$f10.Name = "Segoe UI"
$f10.Size = 14
$f10.Color = 0

# --- Row 11 : new component (MCP3008) ---
$ws.Range("D11").Value = "MCP3008"
$ws.Range("B11").Value = "CAN SPI"
$ws.Range("E11").Value = "https://www.adafruit.com/datasheets/MCP3008.pdf"
$ws.Range("C11").Value = 1

$f11 = $ws.Range("D11").Font
$f11.Name = "Segoe UI"
$f11.Size = 14
$f11.Color = 0

# --- Column D now needs to be wide enough for the longer reference text ---
$ws.Columns.Item(4).AutoFit() | Out-Null

# --- Selection, as left by the author after the edit ---
$ws.Range("E13").Select() | Out-Null

Write-Output "feuille de course updated"
